$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 should match the formatting of the existing
# header cells (bold, bordered, centered) -- copy the format from H1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill data rows 2..36: I = 1 (constant), J = copy of column H value
for ($r = 2; $r -le 36; $r++) {
    $ws.Cells.Item($r, 9).Value2 = 1
    $hval = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 10).Value2 = $hval
}
